# Database tables creation script+
#
# Renames the generic "Order" table/ERD block to "TransactionOrder", and
# disambiguates the four identically-named "Date" fields scattered across
# the CarSale / TireService / Repair / ChipTuning mini-tables into
# distinct, self-describing names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Order" table header + PK field -> "TransactionOrder" table ---
$ws.Range("F1").Value = "TransactionOrder"
$ws.Range("F2").Value = "TransactionOrderID"

# --- CarSale table: generic "Date" field -> "SaleDate" ---
$ws.Range("F10").Value = "SaleDate"

# --- Repair table: "StartDate" field -> "RepairDate" ---
$ws.Range("J16").Value = "RepairDate"

# --- TireService table: generic "Date" field -> "ServiceDate" ---
$ws.Range("F14").Value = "ServiceDate"

# --- ChipTuning table: generic "Date" field -> "TuningDate" ---
$ws.Range("F20").Value = "TuningDate"

# Widen column F to fit the newly lengthened labels (e.g. "TransactionOrderID")
$ws.Columns.Item(6).ColumnWidth = 17.6

# Restore the author's last selection in the sheet view
$ws.Range("F20").Select()
